# Update the player roster table with the new data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New table data (rows 2-18), matching the published diff.
$data = @(
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Bogdan Bogdanovic", "SG,SF", "Atlanta Hawks"),
    @("Pascal Siakam", "SF,PF", "Indiana Pacers"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Julian Champagnie", "SF,PF", "San Antonio Spurs"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Russell Westbrook", "PG", "Denver Nuggets"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jalen Green", "PG,SG", "Houston Rockets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
